$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing student names (column A) ---
# Preserve the order in which new shared strings are first introduced so the
# sharedStrings table layout matches the real editing session.
$ws.Range("K1").Value = "CERTNO"
$ws.Range("A4").Value = "Chinonyelum Ejimuda"
$ws.Range("A2").Value = "Chinedu  Michael"
$ws.Range("L1").Value = "DATE"
$ws.Range("L2").Value = "20th August, 2020"
$ws.Range("K3").Value = "https://roboticscentre.org/       (20082021,009)"
$ws.Range("K2").Value = "https://roboticscentre.org/    (20082021,008)"
$ws.Range("K4").Value = "https://roboticscentre.org/    (20082021,010)"

# --- New hyperlinks for the CERTNO column (K2:K4) ---
$ws.Hyperlinks.Add($ws.Range("K2"), "https://roboticscentre.org/")
$ws.Hyperlinks.Add($ws.Range("K3"), "https://roboticscentre.org/")
$ws.Hyperlinks.Add($ws.Range("K4"), "https://roboticscentre.org/")

# Reapply the shared Hyperlink cell-style so these reuse the existing style
# instead of the fresh one Hyperlinks.Add allocates.
$ws.Range("K2:K4").Style = "Hyperlink"

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 23.833333333333332
$ws.Columns.Item(2).ColumnWidth = 27.5
$ws.Columns.Item(3).ColumnWidth = 9.833333333333334
$ws.Columns.Item(10).ColumnWidth = 12.0
$ws.Columns.Item(11).ColumnWidth = 37.5
$ws.Columns.Item(12).ColumnWidth = 17.666666666666668

# --- Row height for header row ---
$ws.Rows.Item(1).RowHeight = 19.5

# --- Selection ---
$ws.Range("H6").Select()
